$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 9 new rows above the current row 86 (the blank separator row
#    that precedes "Summary"), pushing everything from the old row 86
#    onward down by 9. Inserting one row at a time at the same index is the
#    most reliable way to shift a block of rows in this COM model.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 9; $i++) {
    $ws.Rows.Item(86).Insert()
}

# ---------------------------------------------------------------------------
# 2) Populate the 9 newly inserted rows (86-94) with the new task data.
#    These rows already carry style index 6 (s="6") on columns A:F, copied
#    automatically from the row above (81-85) by the Insert() operation.
# ---------------------------------------------------------------------------
$newTasks = @(
    @(105, "Add all 32 NFL teams master list to lib/constants.ts", "Feature", "ui-dev", "Completed", "lib/constants.ts"),
    @(106, "Update TeamLogo.tsx to use NFL_TEAMS from constants (all 32 teams)", "Enhancement", "ui-dev", "Completed", "components/TeamLogo.tsx"),
    @(107, "Add color palettes for all 32 NFL teams in use-theme.ts", "Enhancement", "ui-dev", "Completed", "hooks/use-theme.ts"),
    @(108, "Design TeamCombobox UX specification", "UX", "ux-expert", "Completed", "N/A (design spec)"),
    @(109, "Build TeamCombobox component with filtered dropdown and logo preview", "Feature", "team-lead-2", "Completed", "components/TeamCombobox.tsx"),
    @(110, "Replace team inputs in game creation form with TeamCombobox", "Enhancement", "team-lead-2", "Completed", "app/page.tsx"),
    @(111, "Replace team inputs in admin page with TeamCombobox", "Enhancement", "ui-dev", "Completed", "app/game/[gameId]/admin/page.tsx"),
    @(112, "Add city field to NFL_TEAMS and city search in TeamCombobox", "Enhancement", "ui-dev", "Completed", "lib/constants.ts, components/TeamCombobox.tsx"),
    @(113, "Prevent duplicate team selection (same team for row and col)", "Bug Fix", "orchestrator", "Completed", "app/page.tsx, app/game/[gameId]/admin/page.tsx")
)

$r = 86
foreach ($task in $newTasks) {
    $ws.Cells.Item($r, 1).Value = $task[0]
    $ws.Cells.Item($r, 2).Value = $task[1]
    $ws.Cells.Item($r, 3).Value = $task[2]
    $ws.Cells.Item($r, 4).Value = $task[3]
    $ws.Cells.Item($r, 5).Value = $task[4]
    $ws.Cells.Item($r, 6).Value = $task[5]
    $r++
}

# ---------------------------------------------------------------------------
# 3) Update the Summary figures that changed because of the 9 new tasks.
#    (The old "Summary" block now starts at row 96, shifted down by 9.)
# ---------------------------------------------------------------------------

# Total Tasks / Completed: 84 -> 93
$ws.Range("B97").Value = 93
$ws.Range("B98").Value = 93

# By Assignee: team-lead-2 13/13 -> 15/15
$ws.Range("B107").Value = "15/15 completed"

# By Assignee: ui-dev 37/37 -> 42/42
$ws.Range("B109").Value = "42/42 completed"

# By Assignee: orchestrator 5/5 -> 6/6 (row keeps style s="7")
$ws.Range("B110").Value = "6/6 completed"

# ---------------------------------------------------------------------------
# 4) Insert one additional row right after "orchestrator" (now row 110) for
#    the brand-new "ux-expert" assignee entry, before the blank separator
#    row that precedes "By Type".
# ---------------------------------------------------------------------------
$ws.Rows.Item(111).Insert()
$ws.Range("A111").Value = "ux-expert"
$ws.Range("B111").Value = "1/1 completed"
# Match the orchestrator row's formatting (style s="7")
$ws.Range("A110:B110").Copy()
$ws.Range("A111:B111").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5) Update the "By Type" counts that changed.
#    (This block now starts at row 113, after the extra ux-expert insert.)
# ---------------------------------------------------------------------------

# Bug Fix: 14 -> 15
$ws.Range("B115").Value = 15

# Enhancement: 11 -> 16
$ws.Range("B118").Value = 16

# Feature: 24 -> 26
$ws.Range("B119").Value = 26

# ---------------------------------------------------------------------------
# 6) Append the new "UX" row at the very end of the By Type block (row 128),
#    styled the same as the "Docs" row above it (s="7").
# ---------------------------------------------------------------------------
$ws.Range("A128").Value = "UX"
$ws.Range("B128").Value = 1
$ws.Range("A127:B127").Copy()
$ws.Range("A128:B128").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 7) Keep the selection on B83 as in the original file.
# ---------------------------------------------------------------------------
$ws.Range("B83").Select()
